$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the sheets that are no longer present in the target workbook
$wb.Worksheets("TINH PHU YEN").Delete()
$wb.Worksheets("TINH TAY NINH").Delete()
$wb.Worksheets("QUAN 2").Delete()
$wb.Worksheets("TINH LAM DONG").Delete()

$ws = $wb.Worksheets("TINH BINH DINH")
$dateSrc = $wb.Worksheets("Sheet1")

# Fill in the new "TINH BINH DINH" data table
$ws.Range("A1").Value = "4101598746 - Công Ty TNHH Chế Biến Lâm Sản Phúc Khang"
$ws.Range("B1").Value = 44261
$ws.Range("C1").Value = "CÔNG TY TNHH CHẾ BIẾN LÂM SẢN PHÚC KHANG"
$ws.Range("D1").Value = "Thôn Hưng Mỹ 1, Xã Cát Hưng, Huyện Phù Cát, Tỉnh Bình Định"
$ws.Range("E1").Value = 982776409
$ws.Range("F1").Value = "Đào Thảo"

$ws.Range("A2").Value = "4101598739 - Công Ty Cổ Phần Daily Feed Việt Nam"
$ws.Range("B2").Value = 44261
$ws.Range("C2").Value = "CÔNG TY CỔ PHẦN DAILY FEED VIỆT NAM"
$ws.Range("D2").Value = "Đội 17, thôn Tư Cung, Xã Phước Thắng, Huyện Tuy Phước, Tỉnh Bình Định"
$ws.Range("E2").Value = 972537148
$ws.Range("F2").Value = "Phạm Đình Tỵ"

$ws.Range("A3").Value = "4101598721 - Công Ty TNHH Tổng Hợp Thương Mại Dịch Vụ Hoàng Phi"
$ws.Range("B3").Value = 44233
$ws.Range("C3").Value = "CÔNG TY TNHH TỔNG HỢP THƯƠNG MẠI DỊCH VỤ HOÀNG PHI"
$ws.Range("D3").Value = "Số 168 Nguyễn Trác, Phường Nhơn Bình, Thành phố Quy Nhơn, Tỉnh Bình Định"
$ws.Range("E3").Value = 937712887
$ws.Range("F3").Value = "Phạm Lê Hoàng Phi"

$ws.Range("A4").Value = "4101598707 - Công Ty TNHH Kinh Doanh Tổng Hợp Phú An"
$ws.Range("B4").Value = 44233
$ws.Range("C4").Value = "CÔNG TY TNHH KINH DOANH TỔNG HỢP PHÚ AN"
$ws.Range("D4").Value = "02 Cần Vương, Phường Nguyễn Văn Cừ, Thành phố Quy Nhơn, Tỉnh Bình Định"
$ws.Range("E4").Value = 977318399
$ws.Range("F4").Value = "Trần Xuân Chí"

$ws.Range("A5").Value = "4101598658 - Công Ty TNHH Gạch Tuy Nen Nhật Đức"
$ws.Range("B5").Value = 44233
$ws.Range("C5").Value = "CÔNG TY TNHH GẠCH TUY NEN NHẬT ĐỨC"
$ws.Range("D5").Value = "Xóm Nam, thôn Lai Nghi, Xã Bình Nghi, Huyện Tây Sơn, Tỉnh Bình Định"
$ws.Range("E5").Value = 399780288
$ws.Range("F5").Value = "Nguyễn Thị Tình"

$ws.Range("A6").Value = "4101598665 - Công Ty TNHH Dịch Vụ Thương Mại Vận Tải Trung Nam"
$ws.Range("B6").Value = 44233
$ws.Range("C6").Value = "CÔNG TY TNHH DỊCH VỤ THƯƠNG MẠI VẬN TẢI TRUNG NAM"
$ws.Range("D6").Value = "Thôn Bình Trị, Xã Mỹ Quang, Huyện Phù Mỹ, Tỉnh Bình Định"
$ws.Range("E6").Value = 965154585
$ws.Range("F6").Value = "Nguyễn Trọng Khiêm"

# Apply the existing short-date style (as used on Sheet1!B24) to the date column
$dateSrc.Range("B24").Copy()
$ws.Range("B1:B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Activate()
